# Updated RAD Test Scripts and Test Data for Existing Liability.
#
# 1. Mark every test row (2-48) as "Execute" = "Y" in column C.
# 2. Rename the PaymentType value "Existing Liability w/Notice Number" to
#    "Existing Liability with Notice/Invoice Number" wherever it appears
#    (column D).
# 3. Update the sheet's selection to match the post-edit state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Flag all data rows to execute.
$ws.Range("C2:C48").Value = "Y"

# 2. Rename the PaymentType label across the sheet.
$ws.Cells.Replace("Existing Liability w/Notice Number", "Existing Liability with Notice/Invoice Number")

# 3. Leave the cursor/selection where the author left it after the edit.
$ws.Range("C15:C48").Select()
$ws.Application.ActiveWindow.ScrollRow = 22

$wb.Save()
